# Generate Report for Handoff
# Adds two new handoff entries (c98675c0-... .md and d51ba409-... .png)
# alongside the existing 8578d384-... .png entry across the Overview,
# zh-cn and de-de sheets, updating the "latest handoff" timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item(1)

$ov.Cells.Item(2,1).Value = "8578d384-50ad-4638-bd7e-56df5605022d.png"
$ov.Cells.Item(2,2).Value = "Ready for handoff"
$ov.Cells.Item(2,3).Value = "Ready for handoff"
$ov.Cells.Item(2,4).Value = "2016-52-21 00:52:02"

$ov.Cells.Item(3,1).Value = "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md"
$ov.Cells.Item(3,2).Value = "Ready for handoff"
$ov.Cells.Item(3,3).Value = "Ready for handoff"
$ov.Cells.Item(3,4).Value = "2016-52-21 00:52:02"

$ov.Cells.Item(4,1).Value = "d51ba409-c477-4e99-8a8e-98b4125e69a0.png"
$ov.Cells.Item(4,2).Value = "Ready for handoff"
$ov.Cells.Item(4,3).Value = "Ready for handoff"
$ov.Cells.Item(4,4).Value = "2016-52-21 00:52:02"

$ov.Range("A2").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/8578d384-50ad-4638-bd7e-56df5605022d.png", "", "", "8578d384-50ad-4638-bd7e-56df5605022d.png") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md", "", "", "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/d51ba409-c477-4e99-8a8e-98b4125e69a0.png", "", "", "d51ba409-c477-4e99-8a8e-98b4125e69a0.png") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item(2)

$zh.Cells.Item(2,1).Value = "8578d384-50ad-4638-bd7e-56df5605022d.png"
$zh.Cells.Item(2,2).Value = ".png"
$zh.Cells.Item(2,3).Value = "Ready for handoff"
$zh.Cells.Item(2,4).Value = "6ce24b83dfed1b9ad8d80f108b16a8a8fdfca54c.png"
$zh.Cells.Item(2,5).Value = "2016-03-21 00:51:59"
$zh.Cells.Item(2,8).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(2,9).Value = "IsDependency"
$zh.Cells.Item(2,10).Value = 'e2e\c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md'

$zh.Cells.Item(3,1).Value = "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md"
$zh.Cells.Item(3,2).Value = ".md"
$zh.Cells.Item(3,3).Value = "Ready for handoff"
$zh.Cells.Item(3,4).Value = "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.cf7b95485007b8e854fde060acec143079ad9888.zh-cn.xlf"
$zh.Cells.Item(3,5).Value = "2016-03-21 00:51:59"
$zh.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(3,9).Value = "Include"

$zh.Cells.Item(4,1).Value = "d51ba409-c477-4e99-8a8e-98b4125e69a0.png"
$zh.Cells.Item(4,2).Value = ".png"
$zh.Cells.Item(4,3).Value = "Ready for handoff"
$zh.Cells.Item(4,4).Value = "f05032a6cbead17ac4c8ff4c4f7f18566e959295.png"
$zh.Cells.Item(4,5).Value = "2016-03-21 00:51:59"
$zh.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(4,9).Value = "IsDependency"
$zh.Cells.Item(4,10).Value = 'e2e\c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md'

$zh.Range("A2").Hyperlinks.Delete()
$zh.Range("B2").Hyperlinks.Delete()
$zh.Range("D2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/8578d384-50ad-4638-bd7e-56df5605022d.png", "", "", "8578d384-50ad-4638-bd7e-56df5605022d.png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/8578d384-50ad-4638-bd7e-56df5605022d.png", "", "", ".png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/467fe1da4fb41d0ab213acfeaf8419856298388d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6ce24b83dfed1b9ad8d80f108b16a8a8fdfca54c.png", "", "", "6ce24b83dfed1b9ad8d80f108b16a8a8fdfca54c.png") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md", "", "", "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md", "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/467fe1da4fb41d0ab213acfeaf8419856298388d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c98675c0-fc5e-4ea6-a817-93231f8cd2dc.cf7b95485007b8e854fde060acec143079ad9888.zh-cn.xlf", "", "", "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.cf7b95485007b8e854fde060acec143079ad9888.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/d51ba409-c477-4e99-8a8e-98b4125e69a0.png", "", "", "d51ba409-c477-4e99-8a8e-98b4125e69a0.png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/d51ba409-c477-4e99-8a8e-98b4125e69a0.png", "", "", ".png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/467fe1da4fb41d0ab213acfeaf8419856298388d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f05032a6cbead17ac4c8ff4c4f7f18566e959295.png", "", "", "f05032a6cbead17ac4c8ff4c4f7f18566e959295.png") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item(3)

$de.Cells.Item(2,1).Value = "8578d384-50ad-4638-bd7e-56df5605022d.png"
$de.Cells.Item(2,2).Value = ".png"
$de.Cells.Item(2,3).Value = "Ready for handoff"
$de.Cells.Item(2,4).Value = "6ce24b83dfed1b9ad8d80f108b16a8a8fdfca54c.png"
$de.Cells.Item(2,5).Value = "2016-03-21 00:52:02"
$de.Cells.Item(2,8).Value = "0001-01-01 00:00:00"
$de.Cells.Item(2,9).Value = "IsDependency"
$de.Cells.Item(2,10).Value = 'e2e\c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md'

$de.Cells.Item(3,1).Value = "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md"
$de.Cells.Item(3,2).Value = ".md"
$de.Cells.Item(3,3).Value = "Ready for handoff"
$de.Cells.Item(3,4).Value = "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.cf7b95485007b8e854fde060acec143079ad9888.de-de.xlf"
$de.Cells.Item(3,5).Value = "2016-03-21 00:52:02"
$de.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$de.Cells.Item(3,9).Value = "Include"

$de.Cells.Item(4,1).Value = "d51ba409-c477-4e99-8a8e-98b4125e69a0.png"
$de.Cells.Item(4,2).Value = ".png"
$de.Cells.Item(4,3).Value = "Ready for handoff"
$de.Cells.Item(4,4).Value = "f05032a6cbead17ac4c8ff4c4f7f18566e959295.png"
$de.Cells.Item(4,5).Value = "2016-03-21 00:52:02"
$de.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$de.Cells.Item(4,9).Value = "IsDependency"
$de.Cells.Item(4,10).Value = 'e2e\c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md'

$de.Range("A2").Hyperlinks.Delete()
$de.Range("B2").Hyperlinks.Delete()
$de.Range("D2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/8578d384-50ad-4638-bd7e-56df5605022d.png", "", "", "8578d384-50ad-4638-bd7e-56df5605022d.png") | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/8578d384-50ad-4638-bd7e-56df5605022d.png", "", "", ".png") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f7b58882bf25ee064504e4abf15ab6b0e1b6f34/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6ce24b83dfed1b9ad8d80f108b16a8a8fdfca54c.png", "", "", "6ce24b83dfed1b9ad8d80f108b16a8a8fdfca54c.png") | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md", "", "", "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/c98675c0-fc5e-4ea6-a817-93231f8cd2dc.md", "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f7b58882bf25ee064504e4abf15ab6b0e1b6f34/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c98675c0-fc5e-4ea6-a817-93231f8cd2dc.cf7b95485007b8e854fde060acec143079ad9888.de-de.xlf", "", "", "c98675c0-fc5e-4ea6-a817-93231f8cd2dc.cf7b95485007b8e854fde060acec143079ad9888.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/d51ba409-c477-4e99-8a8e-98b4125e69a0.png", "", "", "d51ba409-c477-4e99-8a8e-98b4125e69a0.png") | Out-Null
$de.Hyperlinks.Add($de.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/4e0a4c2b56db285315d4ef0b978d18901b185529/e2e/d51ba409-c477-4e99-8a8e-98b4125e69a0.png", "", "", ".png") | Out-Null
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f7b58882bf25ee064504e4abf15ab6b0e1b6f34/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f05032a6cbead17ac4c8ff4c4f7f18566e959295.png", "", "", "f05032a6cbead17ac4c8ff4c4f7f18566e959295.png") | Out-Null
